# Maestro.xlsx — "Add files via upload"
# Inserts a new product row (Salchichas Paladini, barcode 7790079010828)
# as row 11 of the "Artículos" sheet, pushing the existing rows 11-56
# down to 12-57 (dimension grows from A1:P56 to A1:P57).
#
# We deliberately avoid EntireRow/Rows.Insert() here: in this engine it
# synthesizes a brand-new cellXf (merging the row-above's format) for the
# freshly inserted blank row, which would add an unwanted entry to
# styles.xml that the source diff does not contain.
#
# We also avoid a single multi-row Range.Copy(destination) (e.g. copying
# A11:P56 in one shot to A12:P57): this engine's Copy only honours the
# *first* source cell's style and stamps it across the whole destination
# block, which clobbers the s="2"/s="3" style boundary that exists partway
# down this column (rows 2-25 vs 26-56 use different cellXf indices for
# column O). Copying one row at a time — walking from the bottom row
# upward so each source row is read before it gets overwritten — shifts
# both values and their individual per-row styles faithfully.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 56; $r -ge 11; $r--) {
    $src = $ws.Range("A" + $r + ":P" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":P" + ($r + 1))
    $src.Copy($dst)
}

$ws.Range("A11").Value = 7790079010828
$ws.Range("B11").Value = "Salchichas"
$ws.Range("C11").Value = "tipo viena"
$ws.Range("D11").Value = "sin piel"
$ws.Range("E11").Value = "Paladini"
$ws.Range("F11").Value = 6
$ws.Range("G11").Value = "und."
$ws.Range("H11").Value = "bolsa"
$ws.Range("I11").Value = "Salchichas"
$ws.Range("J11").Value = "Argentina"
$ws.Range("K11").Value = 12
$ws.Range("L11").Value = $false
$ws.Range("M11").Value = $true
$ws.Range("N11").Value = "C:\VentaSoft\Imágenes de artículos\7790079010828.png"
$ws.Range("O11").Value = $true
$ws.Range("P11").Value = $true
